$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vals_B = @(7.794976100260179, 7.722960676774888, 7.680252210733521, 7.663246776962808, 7.66044765646529, 7.68002123084283, 7.769843593202482, 7.957090146803507, 8.100228967412585, 8.166293759782128, 8.191425301884765, 8.186008049547164, 8.168359185361794, 8.157562960270949, 8.095928682374129, 8.058345470248854, 8.036819381341308, 8.02954729265697, 8.062337034382409, 8.17354016227512, 8.246873404108861, 8.207681460595898, 8.060532194240029, 7.905368642200784)
$vals_D = @(10.21134977393889, 9.843964667129578, 9.61730752580022, 9.524819305854381, 9.50945877992061, 9.61606049406857, 10.08498537443681, 10.99003944238269, 11.63909723824318, 11.92969303085821, 12.03897819729482, 12.01547660149138, 11.93869961311821, 11.89157053197139, 11.62000459220004, 11.45214764026961, 11.35516613989673, 11.3222580640315, 11.47006194715335, 11.96127205248705, 12.27786150784668, 12.10932384606093, 11.46196437211562, 10.74752780847147)
$vals_E = @(8.731730803213969, 8.354471742924877, 8.113523636775536, 8.013056337279645, 7.996238191207623, 8.112177837637931, 8.603618032096913, 9.49126123316945, 10.09474622687927, 10.35839749140274, 10.45665289777127, 10.43556251615968, 10.36651285100293, 10.32401122918804, 10.07729531647069, 9.92314138825764, 9.83345307773647, 9.802911141614622, 9.939657379104082, 10.38683751570137, 10.66986018560142, 10.51965561808117, 9.932193809614908, 9.25954274566117)
$vals_F = @(76.2614930007023, 75.42697612543434, 74.94744898485799, 74.76037229184604, 74.72981334194499, 74.94489218096223, 75.96697717784203, 78.22813521309502, 80.03960704876262, 80.89467457797231, 81.22277469112397, 81.15192354964776, 80.92158221274724, 80.78104735988445, 79.98433736482352, 79.50341741215246, 79.22973453149083, 79.13757802359952, 79.55431007710918, 80.98912348978622, 81.9518722618955, 81.43579891491648, 79.53129275638926, 77.58953450246631)
$vals_G = @(3.825774641169331, 3.838375321582287, 3.846445018704193, 3.849818092842944, 3.850383325684903, 3.84649016539149, 3.830050835000912, 3.800411650466114, 3.780158746377642, 3.771262132539901, 3.767937581139728, 3.768651626103435, 3.770987736316495, 3.772424420122763, 3.780746436178281, 3.785932015403551, 3.788944486094125, 3.789969619136691, 3.785376918911024, 3.770300368111183, 3.760705239440409, 3.765803087467872, 3.785627780732873, 3.808157859729337)
$vals_J = @(8.621513945549719, 8.638871752506716, 8.650175505366853, 8.654944644573121, 8.655746397046798, 8.650239164158235, 8.627365039747263, 8.587619768590663, 8.561514204464851, 8.550306009926489, 8.546157420516671, 8.547046640662225, 8.549962786216243, 8.551761466541388, 8.562260093278333, 8.568871402609972, 8.572736882345081, 8.574056466308658, 8.568161116712538, 8.549103647997361, 8.537206254743612, 8.543505163302726, 8.568482035888906, 8.597826907568413)
$vals_M = @(38.77136098768339, 38.04842021022784, 37.6115184348304, 37.4354336500292, 37.4063188436326, 37.60913551107266, 38.52076084946198, 40.35514291978796, 41.71982615993626, 42.34201511233875, 42.57763829889134, 42.52689486351262, 42.36140072481703, 42.26002724142523, 41.67917645057818, 41.32306326317989, 41.11837750224893, 41.04910449022122, 41.36095889332749, 42.41001142016135, 43.09562957543309, 42.72976067252819, 41.34382611379048, 39.85510042259699)
$vals_N = @(18.22418924973784, 18.28922054129934, 18.33178553636429, 18.3497977147932, 18.35282899703755, 18.33202575074418, 18.24606737965781, 18.09823932710242, 18.00202628812243, 17.96089437206667, 17.94569336533745, 17.94895057592906, 17.95963629027356, 17.96623027458831, 18.00476703268599, 18.02908018380808, 18.04331282824034, 18.0481745424176, 18.02646633342388, 17.95648750075197, 17.91293459891604, 17.93598136161989, 18.02764726155726, 18.13603588120602)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $vals_B[$i]
    $ws.Range("D$row").Value = $vals_D[$i]
    $ws.Range("E$row").Value = $vals_E[$i]
    $ws.Range("F$row").Value = $vals_F[$i]
    $ws.Range("G$row").Value = $vals_G[$i]
    $ws.Range("J$row").Value = $vals_J[$i]
    $ws.Range("M$row").Value = $vals_M[$i]
    $ws.Range("N$row").Value = $vals_N[$i]
}